$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales Report Sheet")

# Clear out the "July" column (column I) sales figures for the products
# in both tables on the sheet - set them to 0.
$ws.Range("I8").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("I11").Value = 0

$ws.Range("I31").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("I34").Value = 0

# Mirror the selection change captured in the saved file.
$ws.Range("Q29").Select()
